# "Finished the top 10" -- adds the remaining #8..#1 (+ closing blank) slides
# to the "9.5 Things I Love/Hate About Actor Framework" deck, and touches up
# the titles of the two slides (#9.5, #9) that already existed.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Slide 5 ("9.5 No Training Course") -> "#9.5 Love/Hate -No  Training Course"
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$title5 = $s5.Shapes.Item(1)
$tr5 = $title5.TextFrame.TextRange
$tr5.Text = "#9.5 Love/Hate -No"
$tr5.InsertAfter(" ")
$tr5.InsertAfter("Training Course")
$title5.TextFrame.AutoSize = 2

# ---------------------------------------------------------------------
# 2. Slide 6 ("9 AF is Mind Candy!!") -> "#9 Love: AF is Mind Candy!!"
#    Also drop the stray empty trailing paragraph in the body text.
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$title6 = $s6.Shapes.Item(1)
$tr6 = $title6.TextFrame.TextRange
$tr6.Text = "#9 Love: "
$tr6.InsertAfter("AF is Mind Candy!!")

$body6 = $s6.Shapes.Item(2).TextFrame.TextRange
$extraPara = $body6.Paragraphs(4, 1)
$extraPara.Delete()

# ---------------------------------------------------------------------
# 3. Nine new "Title and Content" slides appended after slide 6, finishing
#    out the countdown from #8 down to #1, plus a blank closer.
# ---------------------------------------------------------------------

# #8
$s = $p.Slides.Add(7, 2)
$tr = $s.Shapes.Item(1).TextFrame.TextRange
$tr.Text = "#8"
$tr.InsertAfter(" ")
$tr.InsertAfter("Hate:Heavy")
$tr.InsertAfter("/Bloat")

# #7
$s = $p.Slides.Add(8, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "#7 Love: Think Differently"

# #6
$s = $p.Slides.Add(9, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "#6 Hate:  Breaks Dataflow"

# #5
$s = $p.Slides.Add(10, 2)
$tr = $s.Shapes.Item(1).TextFrame.TextRange
$tr.Text = "#5 Love:  Focus"
$tr.InsertAfter(" on my app, no software design")
$s.Shapes.Item(1).TextFrame.AutoSize = 2

# #4
$s = $p.Slides.Add(11, 2)
$tr = $s.Shapes.Item(1).TextFrame.TextRange
$tr.Text = "#4 Hate: Things"
$tr.InsertAfter(" that should be simple are a kludge")
$s.Shapes.Item(1).TextFrame.AutoSize = 2

# #3
$s = $p.Slides.Add(12, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "#3 Love:  Inspire others"

# #2
$s = $p.Slides.Add(13, 2)
$tr = $s.Shapes.Item(1).TextFrame.TextRange
$tr.Text = "#2 Hate:  Missing design guidelines – How do I look at AF Code and know"
$tr.InsertAfter(" it is well designed")
$s.Shapes.Item(1).TextFrame.AutoSize = 2

# #1
$s = $p.Slides.Add(14, 2)
$tr = $s.Shapes.Item(1).TextFrame.TextRange
$tr.Text = "#1 Love:  "
$tr.InsertAfter("Inherit")
$tr.InsertAfter(" Process")

# Closing blank slide
$s = $p.Slides.Add(15, 2)

Write-Output ("Slide count now: " + $p.Slides.Count)
